$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 1.17
$ws.Range("K2").Value = 5
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.33
$ws.Range("W2").Value = 13
$ws.Range("Z2").Value = 5
$ws.Range("AB2").Value = 29
$ws.Range("AD2").Value = 8.5
$ws.Range("AG2").Value = 67
$ws.Range("AH2").Value = 51

# Row 6
$ws.Range("G6").Value = 1.62
$ws.Range("H6").Value = 4.3
$ws.Range("I6").Value = 4.1
$ws.Range("N6").Value = 1.26
$ws.Range("R6").Value = 1.3
$ws.Range("S6").Value = 2.92
$ws.Range("T6").Value = 17.5
$ws.Range("U6").Value = 15
$ws.Range("V6").Value = 9.25
$ws.Range("W6").Value = 18
$ws.Range("Y6").Value = 14
$ws.Range("Z6").Value = 29
$ws.Range("AA6").Value = 10.75
$ws.Range("AB6").Value = 11.5
$ws.Range("AD6").Value = 26
$ws.Range("AG6").Value = 75
$ws.Range("AH6").Value = 30
$ws.Range("AI6").Value = 24

# Row 11
$ws.Range("O11").Value = 1.72

# Row 12
$ws.Range("G12").Value = 4.1
$ws.Range("N12").Value = 1.9
$ws.Range("O12").Value = 1.9
$ws.Range("R12").Value = 1.8
$ws.Range("S12").Value = 1.91
$ws.Range("U12").Value = 21
$ws.Range("W12").Value = 41
$ws.Range("AD12").Value = 7.5

$wb.Save()
